$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header row (row 1) for the added columns E:H
$ws.Range("E1").Value = "tel"
$ws.Range("F1").Value = "nationalNo"
$ws.Range("G1").Value = "address"

# Set G2 ("Aziziah Street 53") before H1 ("empcode") so that the shared
# string table gets the same append order as the source workbook.
$ws.Range("G2").Value = "Aziziah Street 53"

$ws.Range("H1").Value = "empcode"

# New data row (row 2) for the added columns E:H
$ws.Range("E2").Value = 44556622
$ws.Range("F2").Value = 29067000555
$ws.Range("H2").Value = 1001

# Match the (best-fit) column widths that Excel computed for the new
# "nationalNo" / "address" columns.
$ws.Columns("F").ColumnWidth = 11.166666666666666
$ws.Columns("G").ColumnWidth = 14.833333333333334

# Leave the selection on the newly added H1 cell, matching the saved file.
$ws.Range("H1").Select()
